# Applies the "falta do rotulo titulo" fix:
#  - Slide 1: title becomes the author line; old subtitle text is cleared
#  - Slides 2-6: titles get their slide-number prefix ("2: ", "3: ", ...)
#  - Slide 3: remove the placeholder picture and add a note paragraph instead

$p = $ppt.ActivePresentation

# ---------- Slide 1 ----------
$s1 = $p.Slides.Item(1)

# Title 1: "Introdução à Inteligência Artificial" -> "Autor: Luiz Carlos de Lemos IA"
$s1.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "Autor: Luiz Carlos de Lemos IA"

# Subtitle 2: remove the "- Autor: ..." text, leaving an empty paragraph
$subtitleTr = $s1.Shapes.Item(2).TextFrame.TextRange
if ($subtitleTr.Length -gt 0) {
    $subtitleTr.Characters(1, $subtitleTr.Length).Delete()
}

# ---------- Slide 2 ----------
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "2: O que é Inteligência Artificial?"

# ---------- Slide 3 ----------
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "3: Como funciona a Inteligência Artificial?"

# Append a new top-level paragraph to the content placeholder while
# preserving the existing bullet paragraphs' indent levels.
$contentShape = $s3.Shapes.Item(2)
$contentTr = $contentShape.TextFrame.TextRange
$paraCount = $contentTr.Paragraphs().Count

$existingText = ""
for ($i = 1; $i -le $paraCount; $i++) {
    $existingText = $existingText + $contentTr.Paragraphs($i, 1).Text
}
$contentTr.Text = $existingText + "`r" + "Depois colocar ilustrativa para auxiliar na explicação"

# Restore the indent level (2 = lvl="1") on the paragraphs that had it originally
for ($i = 2; $i -le $paraCount; $i++) {
    $contentTr.Paragraphs($i, 1).IndentLevel = 2
}

# Remove the placeholder picture (it was the 3rd shape on this slide)
$s3.Shapes.Item(3).Delete()

# ---------- Slide 4 ----------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "4: Benefícios e Desafios da Inteligência Artificial"

# ---------- Slide 5 ----------
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "5: Futuro da Inteligência Artificial"

# ---------- Slide 6 ----------
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "6: Conclusão"
